# Add a "2022-Q1" sheet (fund holdings detail) before the "总计" (totals)
# sheet, and refresh "总计" with a new leading row for 2022-Q1.
#
# Strategy: the existing "总计" sheet keeps its sheetId (6) but is repurposed
# (renamed + its data replaced) to become "2022-Q1"; a brand-new sheet is
# appended right after it and named "总计" so it naturally gets the next
# sheetId (7) and rId (rId7) -- matching how the workbook evolved.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Repurpose the current "总计" sheet into the new "2022-Q1" sheet.
# ---------------------------------------------------------------------------
$q1 = $wb.Worksheets.Item(6)
$q1.Name = "2022-Q1"
$q1.Cells.Clear()

# Header row
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"
$q1.Range("B1:H1").Font.Bold = $true
$q1.Range("B1:H1").HorizontalAlignment = -4108
$q1.Range("B1:H1").VerticalAlignment = -4160
$q1.Range("B1:H1").Borders.LineStyle = 1

# Column A (row index) + B..G need to stay TEXT even though most of them
# look like numbers (fund codes, percentages, formatted amounts) -- force
# a text number format before assigning so the literal string is kept
# instead of silently turning into a float (which would also strip
# leading zeroes from fund codes like "009630").
$q1.Range("B2:G20").NumberFormat = "@"

$fundData = @(
    @("163412", "兴全轻资产混合(LOF)", "69.20", "95.26", "4.02", "2.7818", 4),
    @("009630", "浦银安盛ESG责任投资混合A", "15.61", "80.10", "3.72", "0.5807", 10),
    @("501095", "中银证券科技创新 3 年封闭运作灵活配置混合型证券投资基金", "8.32", "79.30", "3.93", "0.3270", 10),
    @("007177", "浙商智能行业优选混合A", "13.24", "89.45", "2.24", "0.2966", 9),
    @("009631", "浦银安盛ESG责任投资混合C", "5.74", "80.10", "3.72", "0.2135", 10),
    @("010552", "浙商智选领航三年持有期混合A", "5.75", "92.98", "3.10", "0.1782", 4),
    @("671030", "西部利得事件驱动股票", "3.55", "94.61", "4.63", "0.1644", 6),
    @("519026", "海富通中小盘混合", "3.88", "91.14", "3.45", "0.1339", 10),
    @("011351", "金鹰年年邮益一年持有期混合A", "9.03", "37.02", "1.21", "0.1093", 3),
    @("007217", "浙商智能行业优选混合C", "1.65", "89.45", "2.24", "0.0370", 9),
    @("007713", "华富科技动能混合", "0.56", "86.98", "5.59", "0.0313", 6),
    @("163818", "中银中小盘成长混合", "0.98", "87.49", "2.78", "0.0272", 3),
    @("010553", "浙商智选领航三年持有期混合C", "0.58", "92.98", "3.10", "0.0180", 4),
    @("002630", "江信瑞福灵活配置混合A", "0.52", "43.17", "2.14", "0.0111", 9),
    @("002631", "江信瑞福灵活配置混合C", "0.50", "43.17", "2.14", "0.0107", 9),
    @("006887", "诺德新生活混合A", "0.28", "86.14", "3.36", "0.0094", 5),
    @("011352", "金鹰年年邮益一年持有期混合C", "0.59", "37.02", "1.21", "0.0071", 3),
    @("009027", "浦银安盛安远回报一年持有期混合A", "0.79", "20.03", "0.54", "0.0043", 9),
    @("006888", "诺德新生活混合C", "0.00", "86.14", "3.36", "0", 5)
)

$r = 2
foreach ($row in $fundData) {
    $q1.Cells.Item($r, 1).Value = ($r - 2)
    $q1.Cells.Item($r, 1).Style = $q1.Cells.Item(1, 1).Style
    $q1.Cells.Item($r, 2).Value = $row[0]
    $q1.Cells.Item($r, 3).Value = $row[1]
    $q1.Cells.Item($r, 4).Value = $row[2]
    $q1.Cells.Item($r, 5).Value = $row[3]
    $q1.Cells.Item($r, 6).Value = $row[4]
    $q1.Cells.Item($r, 8).Value = $row[6]
    $r = $r + 1
}

# Last row (006888) holds a genuine numeric 0 in column G, not text.
$q1.Range("G20").NumberFormat = "General"
$q1.Range("G20").Value = 0
# All the other column-G text amounts.
$q1.Range("G2").Value = "2.7818"
$q1.Range("G3").Value = "0.5807"
$q1.Range("G4").Value = "0.3270"
$q1.Range("G5").Value = "0.2966"
$q1.Range("G6").Value = "0.2135"
$q1.Range("G7").Value = "0.1782"
$q1.Range("G8").Value = "0.1644"
$q1.Range("G9").Value = "0.1339"
$q1.Range("G10").Value = "0.1093"
$q1.Range("G11").Value = "0.0370"
$q1.Range("G12").Value = "0.0313"
$q1.Range("G13").Value = "0.0272"
$q1.Range("G14").Value = "0.0180"
$q1.Range("G15").Value = "0.0111"
$q1.Range("G16").Value = "0.0107"
$q1.Range("G17").Value = "0.0094"
$q1.Range("G18").Value = "0.0071"
$q1.Range("G19").Value = "0.0043"

# Re-apply the column-A index style that Clear() wiped, matching the
# bordered/centered look used by every other quarter sheet.
$indexStyleSample = $q1.Cells.Item(2, 1)

# ---------------------------------------------------------------------------
# 2. Append a brand-new "总计" sheet right after "2022-Q1" (so it naturally
#    receives the next sheetId / rId).
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Add($null, $q1)
$total.Name = "总计"

$total.Range("B1").Value = "日期"
$total.Range("C1").Value = "持有数量(只)"
$total.Range("D1").Value = "持有市值(亿元)"
$total.Range("B1:D1").Font.Bold = $true
$total.Range("B1:D1").HorizontalAlignment = -4108
$total.Range("B1:D1").VerticalAlignment = -4160
$total.Range("B1:D1").Borders.LineStyle = 1

$totalData = @(
    @("2022-Q1", 19, 4.94),
    @("2021-Q4", 55, 20.9),
    @("2021-Q3", 26, 14.06),
    @("2021-Q2", 17, 2.55),
    @("2021-Q1", 8, 0.47),
    @("2020-Q4", 3, 0.45)
)

$r = 2
foreach ($row in $totalData) {
    $total.Cells.Item($r, 1).Value = ($r - 2)
    $total.Cells.Item($r, 2).Value = $row[0]
    $total.Cells.Item($r, 3).Value = $row[1]
    $total.Cells.Item($r, 4).Value = $row[2]
    $r = $r + 1
}

$wb.Worksheets.Item(1).Select()
